$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 84, column A (date serial) ---
$ws.Cells.Item(84, 1).Value = 45457.2916666667

# --- Add new row 85 ---

# A85: date serial, reuse the same number-format style as A84/A83 (no new style created)
$ws.Cells.Item(85, 1).Value = 45460.5883449074
$ws.Cells.Item(84, 1).Copy()
$ws.Cells.Item(85, 1).PasteSpecial(-4122)  # xlPasteFormats

# B85: volume
$ws.Cells.Item(85, 2).Value = 4500

# C85: high
$ws.Cells.Item(85, 3).Value = 6.11999988555908

# D85: low
$ws.Cells.Item(85, 4).Value = 6.01999998092651

# E85: open
$ws.Cells.Item(85, 5).Value = 6.01999998092651

# F85: close
$ws.Cells.Item(85, 6).Value = 6.03999996185303

# G85: adj_close - stored as a shared string of the numeric text (matches source data
# convention elsewhere in the sheet), without introducing a new cell style. Build the
# text via a helper formula cell, then paste-special VALUES only into G85 so the
# General/no-style formatting of the destination is preserved.
$ws.Cells.Item(1, 20).Formula = '=TEXT(6.03999996185303,"0.00000000000000")'
$ws.Cells.Item(1, 20).Copy()
$ws.Cells.Item(85, 7).PasteSpecial(-4163)  # xlPasteValues
$ws.Cells.Item(1, 20).Clear()

# H85: ticker
$ws.Cells.Item(85, 8).Value = "PAL.MI"
